# Updated tijdsbesteding: bump BLOK (B13) from 5 to 10 hours, and add
# 4 + 3 hours to the "week 01-07/02/2016" entry (B14), whose formula
# becomes 7.5+4+3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 10
$ws.Range("B14").Formula = "=7.5+4+3"

# Move the active selection to H15 (matches the diff's <selection activeCell="H15" .../>)
$ws.Range("H15").Select()
